# "Refined metadata to be additional tab"
#
# 1. Update the F-column (time_taken) timestamps on the "data" sheet to a
#    later query time.
# 2. Add a new "metadata" worksheet (after "data") describing the panel
#    query itself (data_name, data_id, data_version, data_version_created,
#    panel_query_time, panel_get_request).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1. refresh the per-row "time_taken" timestamps -----------------------
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:33:15.925251"
$ws.Cells.Item(3, 6).Value = "2021-10-05 14:33:15.925259"
$ws.Cells.Item(4, 6).Value = "2021-10-05 14:33:15.925262"
$ws.Cells.Item(5, 6).Value = "2021-10-05 14:33:15.925264"
$ws.Cells.Item(6, 6).Value = "2021-10-05 14:33:15.925267"
$ws.Cells.Item(7, 6).Value = "2021-10-05 14:33:15.925270"
$ws.Cells.Item(8, 6).Value = "2021-10-05 14:33:15.925272"
$ws.Cells.Item(9, 6).Value = "2021-10-05 14:33:15.925275"
$ws.Cells.Item(10, 6).Value = "2021-10-05 14:33:15.925278"
$ws.Cells.Item(11, 6).Value = "2021-10-05 14:33:15.925280"
$ws.Cells.Item(12, 6).Value = "2021-10-05 14:33:15.925283"
$ws.Cells.Item(13, 6).Value = "2021-10-05 14:33:15.925285"
$ws.Cells.Item(14, 6).Value = "2021-10-05 14:33:15.925288"
$ws.Cells.Item(15, 6).Value = "2021-10-05 14:33:15.925290"
$ws.Cells.Item(16, 6).Value = "2021-10-05 14:33:15.925293"
$ws.Cells.Item(17, 6).Value = "2021-10-05 14:33:15.925295"
$ws.Cells.Item(18, 6).Value = "2021-10-05 14:33:15.925298"
$ws.Cells.Item(19, 6).Value = "2021-10-05 14:33:15.925301"
$ws.Cells.Item(20, 6).Value = "2021-10-05 14:33:15.925303"
$ws.Cells.Item(21, 6).Value = "2021-10-05 14:33:15.925306"
$ws.Cells.Item(22, 6).Value = "2021-10-05 14:33:15.925308"
$ws.Cells.Item(23, 6).Value = "2021-10-05 14:33:15.925311"
$ws.Cells.Item(24, 6).Value = "2021-10-05 14:33:15.925313"
$ws.Cells.Item(25, 6).Value = "2021-10-05 14:33:15.925316"
$ws.Cells.Item(26, 6).Value = "2021-10-05 14:33:15.925319"
$ws.Cells.Item(27, 6).Value = "2021-10-05 14:33:15.925321"

# --- 2. add the "metadata" sheet right after "data" ------------------------
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# header row formatting (bold, centered, thin border) -- matches the
# formatting already used for the "data" sheet's header row
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$meta.Cells.Item(2, 1).Value = 0
$meta.Cells.Item(2, 1).Font.Bold = $true
$meta.Cells.Item(2, 1).HorizontalAlignment = -4108
$meta.Cells.Item(2, 1).VerticalAlignment = -4160
$meta.Cells.Item(2, 1).Borders.LineStyle = 1

$meta.Cells.Item(2, 2).Value = "Bardet Biedl syndrome"
$meta.Cells.Item(2, 3).Value = 53

# keep "1.10" as literal text (not coerced to the number 1.1)
$meta.Cells.Item(2, 4).NumberFormat = "@"
$meta.Cells.Item(2, 4).Value = "1.10"

$meta.Cells.Item(2, 5).Value = "2021-08-06T01:14:15.525045Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:33:15.922010"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/53/?format=json"

$ws.Select()
